$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 298: data now matches what was previously row 300
$ws.Range("B298").Value = 7083188
$row298 = New-Object 'object[,]' 1,26
$row298[0,0] = 'Legia Warsaw'
$row298[0,1] = 'Zaglebie Lubin'
$row298[0,2] = 2
$row298[0,3] = 1
$row298[0,4] = 2
$row298[0,5] = 0
$row298[0,6] = 'H'
$row298[0,7] = 1.5
$row298[0,8] = 4
$row298[0,9] = 5.5
$row298[0,10] = 1.6
$row298[0,11] = 4.1
$row298[0,12] = 4.333
$row298[0,13] = -0.75
$row298[0,14] = 1.825
$row298[0,15] = 2.025
$row298[0,16] = 3
$row298[0,17] = 1.875
$row298[0,18] = 1.975
$row298[0,19] = 0.6000000000000001
$row298[0,20] = -1
$row298[0,21] = -1
$row298[0,22] = 0.4125
$row298[0,23] = -0.5
$row298[0,24] = 0
$row298[0,25] = 0
$ws.Range("E298:AD298").Value = $row298

# Row 299: data now matches what was previously row 301
$ws.Range("B299").Value = 7093820
$row299 = New-Object 'object[,]' 1,26
$row299[0,0] = 'Ruch Chorzow'
$row299[0,1] = 'Cracovia Krakow'
$row299[0,2] = 2
$row299[0,3] = 0
$row299[0,4] = 1
$row299[0,5] = 0
$row299[0,6] = 'H'
$row299[0,7] = 2.5
$row299[0,8] = 3.4
$row299[0,9] = 2.5
$row299[0,10] = 2.6
$row299[0,11] = 3.6
$row299[0,12] = 2.3
$row299[0,13] = 0
$row299[0,14] = 2.025
$row299[0,15] = 1.825
$row299[0,16] = 3
$row299[0,17] = 2.025
$row299[0,18] = 1.825
$row299[0,19] = 1.6
$row299[0,20] = -1
$row299[0,21] = -1
$row299[0,22] = 1.025
$row299[0,23] = -1
$row299[0,24] = -1
$row299[0,25] = 0.825
$ws.Range("E299:AD299").Value = $row299

# Row 300: data now matches what was previously row 299
$ws.Range("B300").Value = 7083187
$row300 = New-Object 'object[,]' 1,26
$row300[0,0] = 'Lech Poznan'
$row300[0,1] = 'Korona Kielce'
$row300[0,2] = 1
$row300[0,3] = 2
$row300[0,4] = 1
$row300[0,5] = 0
$row300[0,6] = 'A'
$row300[0,7] = 1.8
$row300[0,8] = 3.8
$row300[0,9] = 3.6
$row300[0,10] = 2.1
$row300[0,11] = 3.7
$row300[0,12] = 2.9
$row300[0,13] = -0.25
$row300[0,14] = 1.9
$row300[0,15] = 1.95
$row300[0,16] = 2.75
$row300[0,17] = 1.925
$row300[0,18] = 1.925
$row300[0,19] = -1
$row300[0,20] = -1
$row300[0,21] = 1.9
$row300[0,22] = -1
$row300[0,23] = 0.95
$row300[0,24] = 0.4625
$row300[0,25] = -0.5
$ws.Range("E300:AD300").Value = $row300

# Row 301: data now matches what was previously row 298
$ws.Range("B301").Value = 7041338
$row301 = New-Object 'object[,]' 1,26
$row301[0,0] = 'Jagiellonia Bialystok'
$row301[0,1] = 'Warta Poznan'
$row301[0,2] = 3
$row301[0,3] = 0
$row301[0,4] = 3
$row301[0,5] = 0
$row301[0,6] = 'H'
$row301[0,7] = 1.444
$row301[0,8] = 4.75
$row301[0,9] = 5.25
$row301[0,10] = 1.4
$row301[0,11] = 4.75
$row301[0,12] = 5.75
$row301[0,13] = -1.25
$row301[0,14] = 1.9
$row301[0,15] = 1.95
$row301[0,16] = 3
$row301[0,17] = 1.925
$row301[0,18] = 1.925
$row301[0,19] = 0.3999999999999999
$row301[0,20] = -1
$row301[0,21] = -1
$row301[0,22] = 0.8999999999999999
$row301[0,23] = -1
$row301[0,24] = 0
$row301[0,25] = 0
$ws.Range("E301:AD301").Value = $row301

# Row 305: data now matches what was previously row 306
$ws.Range("B305").Value = 7074364
$row305 = New-Object 'object[,]' 1,26
$row305[0,0] = 'Rakow Czestochowa'
$row305[0,1] = 'Slask Wroclaw'
$row305[0,2] = 1
$row305[0,3] = 2
$row305[0,4] = 1
$row305[0,5] = 0
$row305[0,6] = 'A'
$row305[0,7] = 2.5
$row305[0,8] = 3.6
$row305[0,9] = 2.4
$row305[0,10] = 2.15
$row305[0,11] = 3.6
$row305[0,12] = 2.875
$row305[0,13] = -0.25
$row305[0,14] = 1.95
$row305[0,15] = 1.9
$row305[0,16] = 2.5
$row305[0,17] = 1.875
$row305[0,18] = 1.975
$row305[0,19] = -1
$row305[0,20] = -1
$row305[0,21] = 1.875
$row305[0,22] = -1
$row305[0,23] = 0.8999999999999999
$row305[0,24] = 0.875
$row305[0,25] = -1
$ws.Range("E305:AD305").Value = $row305

# Row 306: data now matches what was previously row 305
$ws.Range("B306").Value = 7093821
$row306 = New-Object 'object[,]' 1,26
$row306[0,0] = 'LKS Lodz'
$row306[0,1] = 'Stal Mielec'
$row306[0,2] = 3
$row306[0,3] = 2
$row306[0,4] = 3
$row306[0,5] = 0
$row306[0,6] = 'H'
$row306[0,7] = 2.5
$row306[0,8] = 3.4
$row306[0,9] = 2.5
$row306[0,10] = 2.2
$row306[0,11] = 3.5
$row306[0,12] = 2.8
$row306[0,13] = -0.25
$row306[0,14] = 2.025
$row306[0,15] = 1.825
$row306[0,16] = 3
$row306[0,17] = 2
$row306[0,18] = 1.85
$row306[0,19] = 1.2
$row306[0,20] = -1
$row306[0,21] = -1
$row306[0,22] = 1.025
$row306[0,23] = -1
$row306[0,24] = 1
$row306[0,25] = -1
$ws.Range("E306:AD306").Value = $row306
